# Add files via upload
# Adds two new columns (G: Core, H: Elective) to the course list sheet,
# marking select course rows with "Y" in the appropriate column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Core"
$ws.Range("H1").Value = "Elective"

# Match the formatting already used by the other header cells (F1 uses style index 4,
# which carries the bordered/centered "header" look).
$ws.Range("G1").Style = $ws.Range("F1").Style
$ws.Range("H1").Style = $ws.Range("F1").Style

# Mark individual course rows as Core ("G") or Elective ("H") with "Y"
$ws.Range("H8").Value = "Y"
$ws.Range("G9").Value = "Y"
$ws.Range("H10").Value = "Y"
$ws.Range("G11").Value = "Y"
$ws.Range("G15").Value = "Y"
$ws.Range("G18").Value = "Y"
$ws.Range("H19").Value = "Y"
$ws.Range("H20").Value = "Y"
$ws.Range("H22").Value = "Y"
$ws.Range("H23").Value = "Y"
$ws.Range("G24").Value = "Y"
$ws.Range("G25").Value = "Y"
$ws.Range("H26").Value = "Y"
$ws.Range("H27").Value = "Y"
$ws.Range("H28").Value = "Y"
$ws.Range("G30").Value = "Y"
$ws.Range("H33").Value = "Y"
$ws.Range("H34").Value = "Y"
$ws.Range("H36").Value = "Y"
$ws.Range("H37").Value = "Y"
$ws.Range("H38").Value = "Y"
$ws.Range("H40").Value = "Y"
$ws.Range("H45").Value = "Y"
$ws.Range("H51").Value = "Y"
$ws.Range("H53").Value = "Y"
$ws.Range("H54").Value = "Y"
$ws.Range("H55").Value = "Y"
$ws.Range("H56").Value = "Y"
$ws.Range("H58").Value = "Y"
$ws.Range("H59").Value = "Y"
$ws.Range("H61").Value = "Y"
$ws.Range("H62").Value = "Y"
$ws.Range("H63").Value = "Y"
$ws.Range("H64").Value = "Y"
$ws.Range("H65").Value = "Y"
$ws.Range("H66").Value = "Y"
$ws.Range("H68").Value = "Y"
$ws.Range("H69").Value = "Y"
$ws.Range("H71").Value = "Y"
$ws.Range("H73").Value = "Y"

# Reset the sheet view: drop the old "scrolled down to row 124" state and
# move the active selection to H144 (just past the new last data row).
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("H144").Select

# Shrink the saved window size/position to match the authored workbook view.
$excel.ActiveWindow.WindowState = -4143
$excel.Left = -110
$excel.Top = -110
$excel.Width = 19420
$excel.Height = 10420
